$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 76 (pushing existing rows 76-165 down to 78-167)
$ws.Rows.Item(76).Insert()
$ws.Rows.Item(76).Insert()

# Populate new row 76 with the new weekly price record
$ws.Cells.Item(76,1).Value = 11
$ws.Cells.Item(76,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(76,3).Value = "Bíobío"
$ws.Cells.Item(76,4).Value = 44741
$ws.Cells.Item(76,5).Value = 8
$ws.Cells.Item(76,6).Value = "Fruta"
$ws.Cells.Item(76,7).Value = 100101
$ws.Cells.Item(76,8).Value = "Berries"
$ws.Cells.Item(76,9).Value = 100101007
$ws.Cells.Item(76,10).Value = "Kiwi"
$ws.Cells.Item(76,11).Value = "Hayward"
$ws.Cells.Item(76,12).Value = "Primera"
$ws.Cells.Item(76,13).Value = 50
$ws.Cells.Item(76,14).Value = 7000
$ws.Cells.Item(76,15).Value = 7000
$ws.Cells.Item(76,16).Value = 7000
$ws.Cells.Item(76,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(76,18).Value = "Región de O'Higgins"
$ws.Cells.Item(76,19).Value = 389
$ws.Cells.Item(76,20).Value = 18

# Populate new row 77 with the new weekly price record
$ws.Cells.Item(77,1).Value = 11
$ws.Cells.Item(77,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(77,3).Value = "Bíobío"
$ws.Cells.Item(77,4).Value = 44741
$ws.Cells.Item(77,5).Value = 8
$ws.Cells.Item(77,6).Value = "Fruta"
$ws.Cells.Item(77,7).Value = 100101
$ws.Cells.Item(77,8).Value = "Berries"
$ws.Cells.Item(77,9).Value = 100101007
$ws.Cells.Item(77,10).Value = "Kiwi"
$ws.Cells.Item(77,11).Value = "Hayward"
$ws.Cells.Item(77,12).Value = "Segunda"
$ws.Cells.Item(77,13).Value = 50
$ws.Cells.Item(77,14).Value = 6000
$ws.Cells.Item(77,15).Value = 6000
$ws.Cells.Item(77,16).Value = 6000
$ws.Cells.Item(77,17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(77,18).Value = "Región de O'Higgins"
$ws.Cells.Item(77,19).Value = 333
$ws.Cells.Item(77,20).Value = 18
